$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.092947483062744
$ws.Range("B1").Value = 2.307811498641968
$ws.Range("C1").Value = 9.660909652709961
$ws.Range("D1").Value = 2.257486581802368
$ws.Range("E1").Value = 1.294302701950073
